$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Grow "Table3" by one row (5 -> 6 data rows) to make room for the new
#    "Hot Dog" menu item, then rebuild all the data rows with the simplified
#    icon/ingredient text from the commit.
# ---------------------------------------------------------------------------

$lo = $ws.ListObjects.Item(1)
if ($lo.ListRows.Count -lt 5) {
    $lo.ListRows.Add() | Out-Null
}

$data = @(
    @("Cheese Burger", "Sesame Burger Bun / Beef Patty[Beef, water, toasted wheat crumbs, soy protein, flavour, salt, spices, spice extracts, onion, garlic] / Pickles / Cheddar Cheese / Red Onion", "Soy, sesame, wheat, milk, eggs, gluten. Bun made in a facility that also processes: crustacean, shellfish, fish, sulphites, eggs, mustard. ", "Island City Seeded Bun", "BC", "Beef_Burger"),
    @("Veggie Burger", "Sesame Burger Bun / Spicy Black Bean Patty[Textured Soy Protein, Cooked Black Beans, Cooked Brown Rice, Onion, Water, Corn Oil, Egg Whites, Diced Tomatoes, Roasted Corn, Bulgur Wheat, Green Chiles, Modified Milk Ingredients, Spices, Yeast Extract, Natural and Artificial Flavours, Lactic Acid, Calcium Phosphate] / Pickles / Cheddar Cheese / Red Onion", "Soy, wheat, milk, eggs, gluten. Bun made in a facility that also processes: crustacean, shellfish, fish, sulphites, eggs, mustard. ", "Island City Seeded Bun", "BC,VEG", "Veggie_Burger"),
    @("Mac & Cheese Bites", "Falafel / Feta / Red Onion / Cucumber / Mixed Greens / Tomato / Tzatziki", "Wheat, sulphites.", "NA", "NA", "placeholder"),
    @("Mozza Sticks", "Turkey / Bacon / Cheddar Cheese / Lettuce / Tomato / Garlic Aioli", "Wheat, milk.", "NA", "VEG", "placeholder"),
    @("Hot Dog", "Bun / Chicken and Beef weiner", "Wheat.  Bun made in a facility that also processes: crustacean, shellfish, fish, sulphites, eggs, mustard. ", "Island City Hot Dog Bun", "BC, DF", "Hot_Dog")
)

$rowNum = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowNum, 1).Value = $rec[0]
    $ws.Cells.Item($rowNum, 2).Value = $rec[1]
    $ws.Cells.Item($rowNum, 3).Value = $rec[2]
    $ws.Cells.Item($rowNum, 4).Value = $rec[3]
    $ws.Cells.Item($rowNum, 5).Value = $rec[4]
    $ws.Cells.Item($rowNum, 6).Value = $rec[5]
    $rowNum++
}

# ---------------------------------------------------------------------------
# 2. Styling: the old wrapped-text treatment on B3/C3 is gone (back to the
#    default "Normal" style); the ItemName column (A2, A3, A6) now gets a
#    centred (both axes), black, size-11 Calibri look instead.
# ---------------------------------------------------------------------------

$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Style = "Normal"

foreach ($addr in @("A2", "A3", "A6")) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.Font.Size = 11
    $r.Font.Color = 0
}

# Drop the custom row heights that used to accommodate the wrapped text.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# ---------------------------------------------------------------------------
# 3. Selection follows the edit point in the author's session.
# ---------------------------------------------------------------------------

$ws.Range("C5").Select() | Out-Null
